$wb = $excel.ActiveWorkbook

# ---------------- Sheet1 (training) ----------------
$ws1 = $wb.Worksheets.Item(1)

# Step 1: Set D1 header and D2:D6 formulas (percentage, referring to new C column)
$ws1.Range("D1").Value = $ws1.Range("C1").Value2

$ws1.Range("D2:D6").NumberFormat = "0.0%"
$ws1.Range("D2").Formula = "=C2/SUM(C2:C6)"
$ws1.Range("D3").Formula = "=C3/SUM(C2:C6)"
$ws1.Range("D4").Formula = "=C4/SUM(C2:C6)"
$ws1.Range("D5").Formula = "=C5/SUM(C2:C6)"
$ws1.Range("D6").Formula = "=C6/SUM(C2:C6)"

# Step 2: Move the Count values from B to C (plain numbers/text, no formatting)
$ws1.Range("C1:C6").ClearFormats()
$ws1.Range("C1").Value = $ws1.Range("B1").Value2
$ws1.Range("C2").Value = $ws1.Range("B2").Value2
$ws1.Range("C3").Value = $ws1.Range("B3").Value2
$ws1.Range("C4").Value = $ws1.Range("B4").Value2
$ws1.Range("C5").Value = $ws1.Range("B5").Value2
$ws1.Range("C6").Value = $ws1.Range("B6").Value2

# Step 3: Set B to the new Label column (order chosen to match shared-string insertion order)
$ws1.Range("B1").Value = "Label"
$ws1.Range("B2").Value = "Normal"
$ws1.Range("B3").Value = "Mild Non-Proliferative DR"
$ws1.Range("B6").Value = "Proliferative DR"
$ws1.Range("B4").Value = "Moderate Non-Proliferative DR"
$ws1.Range("B5").Value = "Severe Non-Proliferative DR"

$ws1.Columns.Item(2).ColumnWidth = 26.88671875

$ws1.Range("B1:B6").Select()

# ---------------- Sheet2 (test) ----------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("D1").Value = $ws2.Range("C1").Value2

$ws2.Range("D2:D6").NumberFormat = "0.0%"
$ws2.Range("D2").Formula = "=C2/SUM(C2:C6)"
$ws2.Range("D3").Formula = "=C3/SUM(C2:C6)"
$ws2.Range("D4").Formula = "=C4/SUM(C2:C6)"
$ws2.Range("D5").Formula = "=C5/SUM(C2:C6)"
$ws2.Range("D6").Formula = "=C6/SUM(C2:C6)"

$ws2.Range("C1:C6").ClearFormats()
$ws2.Range("C1").Value = $ws2.Range("B1").Value2
$ws2.Range("C2").Value = $ws2.Range("B2").Value2
$ws2.Range("C3").Value = $ws2.Range("B3").Value2
$ws2.Range("C4").Value = $ws2.Range("B4").Value2
$ws2.Range("C5").Value = $ws2.Range("B5").Value2
$ws2.Range("C6").Value = $ws2.Range("B6").Value2

$ws2.Range("B1").Value = "Label"
$ws2.Range("B2").Value = "Normal"
$ws2.Range("B3").Value = "Mild Non-Proliferative DR"
$ws2.Range("B6").Value = "Proliferative DR"
$ws2.Range("B4").Value = "Moderate Non-Proliferative DR"
$ws2.Range("B5").Value = "Severe Non-Proliferative DR"

$ws2.Columns.Item(2).ColumnWidth = 26.88671875

$ws2.Range("B1:B6").Select()

# ---------------- Workbook-level view state ----------------
$ws2.Activate()
